# Rewrites the negative/positive anchor-word confidence tables with the
# "min 5" toy-spam run (one extra row of data, several words replaced).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(1,1).Value = "negative"
$ws.Cells.Item(1,10).Value = "positive"

$ws.Cells.Item(3,1).Value = "poorly"
$ws.Cells.Item(3,2).Value = 0.9782608695652174
$ws.Cells.Item(3,3).Value = 45
$ws.Cells.Item(3,4).Value = 45
$ws.Cells.Item(3,5).Value = 0
$ws.Cells.Item(3,6).Value = 1
$ws.Cells.Item(3,7).Value = $false
$ws.Cells.Item(3,8).Value = 1

$ws.Cells.Item(4,1).Value = "disappointing"
$ws.Cells.Item(4,2).Value = 0.8636363636363636
$ws.Cells.Item(4,3).Value = 38
$ws.Cells.Item(4,4).Value = 38
$ws.Cells.Item(4,5).Value = 0
$ws.Cells.Item(4,6).Value = 1
$ws.Cells.Item(4,7).Value = $false
$ws.Cells.Item(4,8).Value = 6

$ws.Cells.Item(5,1).Value = "disappointed"
$ws.Cells.Item(5,2).Value = 0.7688172043010753
$ws.Cells.Item(5,3).Value = 143
$ws.Cells.Item(5,4).Value = 143
$ws.Cells.Item(5,5).Value = 0
$ws.Cells.Item(5,6).Value = 1
$ws.Cells.Item(5,7).Value = $false
$ws.Cells.Item(5,8).Value = 43

$ws.Cells.Item(6,1).Value = "however"
$ws.Cells.Item(6,2).Value = 0.765625
$ws.Cells.Item(6,3).Value = 49
$ws.Cells.Item(6,4).Value = 49
$ws.Cells.Item(6,5).Value = 0
$ws.Cells.Item(6,6).Value = 1
$ws.Cells.Item(6,7).Value = $false
$ws.Cells.Item(6,8).Value = 15

$ws.Cells.Item(7,1).Value = "broke"
$ws.Cells.Item(7,2).Value = 0.7427184466019418
$ws.Cells.Item(7,3).Value = 153
$ws.Cells.Item(7,4).Value = 153
$ws.Cells.Item(7,5).Value = 0
$ws.Cells.Item(7,6).Value = 1
$ws.Cells.Item(7,7).Value = $false
$ws.Cells.Item(7,8).Value = 53

$ws.Cells.Item(8,1).Value = "poor"
$ws.Cells.Item(8,2).Value = 0.7323943661971831
$ws.Cells.Item(8,3).Value = 52
$ws.Cells.Item(8,4).Value = 52
$ws.Cells.Item(8,5).Value = 0
$ws.Cells.Item(8,6).Value = 1
$ws.Cells.Item(8,7).Value = $false
$ws.Cells.Item(8,8).Value = 19

$ws.Cells.Item(9,1).Value = "junk"
$ws.Cells.Item(9,2).Value = 0.6909090909090909
$ws.Cells.Item(9,3).Value = 38
$ws.Cells.Item(9,4).Value = 38
$ws.Cells.Item(9,5).Value = 0
$ws.Cells.Item(9,6).Value = 1
$ws.Cells.Item(9,7).Value = $false
$ws.Cells.Item(9,8).Value = 17

$ws.Cells.Item(10,1).Value = "waste"
$ws.Cells.Item(10,2).Value = 0.6621621621621622
$ws.Cells.Item(10,3).Value = 98
$ws.Cells.Item(10,4).Value = 98
$ws.Cells.Item(10,5).Value = 0
$ws.Cells.Item(10,6).Value = 1
$ws.Cells.Item(10,7).Value = $false
$ws.Cells.Item(10,8).Value = 50

$ws.Cells.Item(11,1).Value = "smaller"
$ws.Cells.Item(11,2).Value = 0.5630252100840336
$ws.Cells.Item(11,3).Value = 67
$ws.Cells.Item(11,4).Value = 67
$ws.Cells.Item(11,5).Value = 0
$ws.Cells.Item(11,6).Value = 1
$ws.Cells.Item(11,7).Value = $false
$ws.Cells.Item(11,8).Value = 52

$ws.Cells.Item(12,1).Value = "broken"
$ws.Cells.Item(12,2).Value = 0.5180722891566265
$ws.Cells.Item(12,3).Value = 43
$ws.Cells.Item(12,4).Value = 43
$ws.Cells.Item(12,5).Value = 0
$ws.Cells.Item(12,6).Value = 1
$ws.Cells.Item(12,7).Value = $false
$ws.Cells.Item(12,8).Value = 40

$ws.Cells.Item(13,1).Value = "small"
$ws.Cells.Item(13,2).Value = 0.5072463768115942
$ws.Cells.Item(13,3).Value = 175
$ws.Cells.Item(13,4).Value = 175
$ws.Cells.Item(13,5).Value = 0
$ws.Cells.Item(13,6).Value = 1
$ws.Cells.Item(13,7).Value = $false
$ws.Cells.Item(13,8).Value = 170

$ws.Cells.Item(14,1).Value = "apart"
$ws.Cells.Item(14,2).Value = 0.4947368421052631
$ws.Cells.Item(14,3).Value = 47
$ws.Cells.Item(14,4).Value = 47
$ws.Cells.Item(14,5).Value = 0
$ws.Cells.Item(14,6).Value = 1
$ws.Cells.Item(14,7).Value = $false
$ws.Cells.Item(14,8).Value = 48

$ws.Cells.Item(15,1).Value = "plastic"
$ws.Cells.Item(15,2).Value = 0.4645669291338583
$ws.Cells.Item(15,3).Value = 59
$ws.Cells.Item(15,4).Value = 59
$ws.Cells.Item(15,5).Value = 0
$ws.Cells.Item(15,6).Value = 1
$ws.Cells.Item(15,7).Value = $false
$ws.Cells.Item(15,8).Value = 68

$ws.Cells.Item(16,1).Value = "difficult"
$ws.Cells.Item(16,2).Value = 0.3820224719101123
$ws.Cells.Item(16,3).Value = 34
$ws.Cells.Item(16,4).Value = 34
$ws.Cells.Item(16,5).Value = 0
$ws.Cells.Item(16,6).Value = 1
$ws.Cells.Item(16,7).Value = $false
$ws.Cells.Item(16,8).Value = 55

$ws.Cells.Item(17,1).Value = "thought"
$ws.Cells.Item(17,2).Value = 0.3415841584158416
$ws.Cells.Item(17,3).Value = 69
$ws.Cells.Item(17,4).Value = 69
$ws.Cells.Item(17,5).Value = 0
$ws.Cells.Item(17,6).Value = 1
$ws.Cells.Item(17,7).Value = $false
$ws.Cells.Item(17,8).Value = 133

$ws.Cells.Item(18,1).Value = "ok"
$ws.Cells.Item(18,2).Value = 0.296875
$ws.Cells.Item(18,3).Value = 38
$ws.Cells.Item(18,4).Value = 38
$ws.Cells.Item(18,5).Value = 0
$ws.Cells.Item(18,6).Value = 1
$ws.Cells.Item(18,7).Value = $false
$ws.Cells.Item(18,8).Value = 90

$ws.Cells.Item(19,1).Value = "cheap"
$ws.Cells.Item(19,2).Value = 0.2748815165876777
$ws.Cells.Item(19,3).Value = 58
$ws.Cells.Item(19,4).Value = 58
$ws.Cells.Item(19,5).Value = 0
$ws.Cells.Item(19,6).Value = 1
$ws.Cells.Item(19,7).Value = $false
$ws.Cells.Item(19,8).Value = 153

$ws.Cells.Item(20,1).Value = "size"
$ws.Cells.Item(20,2).Value = 0.2680412371134021
$ws.Cells.Item(20,3).Value = 52
$ws.Cells.Item(20,4).Value = 52
$ws.Cells.Item(20,5).Value = 0
$ws.Cells.Item(20,6).Value = 1
$ws.Cells.Item(20,7).Value = $false
$ws.Cells.Item(20,8).Value = 142

$ws.Cells.Item(21,1).Value = "though"
$ws.Cells.Item(21,2).Value = 0.2478632478632479
$ws.Cells.Item(21,3).Value = 29
$ws.Cells.Item(21,4).Value = 29
$ws.Cells.Item(21,5).Value = 0
$ws.Cells.Item(21,6).Value = 1
$ws.Cells.Item(21,7).Value = $false
$ws.Cells.Item(21,8).Value = 88

$ws.Cells.Item(22,1).Value = "hard"
$ws.Cells.Item(22,2).Value = 0.215
$ws.Cells.Item(22,3).Value = 43
$ws.Cells.Item(22,4).Value = 43
$ws.Cells.Item(22,5).Value = 0
$ws.Cells.Item(22,6).Value = 1
$ws.Cells.Item(22,7).Value = $false
$ws.Cells.Item(22,8).Value = 157

$ws.Cells.Item(23,1).Value = "item"
$ws.Cells.Item(23,2).Value = 0.2065217391304348
$ws.Cells.Item(23,3).Value = 57
$ws.Cells.Item(23,4).Value = 57
$ws.Cells.Item(23,5).Value = 0
$ws.Cells.Item(23,6).Value = 1
$ws.Cells.Item(23,7).Value = $false
$ws.Cells.Item(23,8).Value = 219

$ws.Cells.Item(24,1).Value = "money"
$ws.Cells.Item(24,2).Value = 0.1930379746835443
$ws.Cells.Item(24,3).Value = 61
$ws.Cells.Item(24,4).Value = 61
$ws.Cells.Item(24,5).Value = 0
$ws.Cells.Item(24,6).Value = 1
$ws.Cells.Item(24,7).Value = $false
$ws.Cells.Item(24,8).Value = 255

$ws.Cells.Item(25,1).Value = "would"
$ws.Cells.Item(25,2).Value = 0.1928783382789317
$ws.Cells.Item(25,3).Value = 130
$ws.Cells.Item(25,4).Value = 130
$ws.Cells.Item(25,5).Value = 0
$ws.Cells.Item(25,6).Value = 1
$ws.Cells.Item(25,7).Value = $false
$ws.Cells.Item(25,8).Value = 544

$ws.Cells.Item(26,1).Value = "work"
$ws.Cells.Item(26,2).Value = 0.1708860759493671
$ws.Cells.Item(26,3).Value = 54
$ws.Cells.Item(26,4).Value = 54
$ws.Cells.Item(26,5).Value = 0
$ws.Cells.Item(26,6).Value = 1
$ws.Cells.Item(26,7).Value = $false
$ws.Cells.Item(26,8).Value = 262

$ws.Cells.Item(27,1).Value = "better"
$ws.Cells.Item(27,2).Value = 0.1588785046728972
$ws.Cells.Item(27,3).Value = 34
$ws.Cells.Item(27,4).Value = 34
$ws.Cells.Item(27,5).Value = 0
$ws.Cells.Item(27,6).Value = 1
$ws.Cells.Item(27,7).Value = $false
$ws.Cells.Item(27,8).Value = 180

$ws.Cells.Item(28,1).Value = "product"
$ws.Cells.Item(28,2).Value = 0.1563876651982379
$ws.Cells.Item(28,3).Value = 71
$ws.Cells.Item(28,4).Value = 71
$ws.Cells.Item(28,5).Value = 0
$ws.Cells.Item(28,6).Value = 1
$ws.Cells.Item(28,7).Value = $false
$ws.Cells.Item(28,8).Value = 383

$ws.Cells.Item(29,1).Value = "3"
$ws.Cells.Item(29,2).Value = 0.1411290322580645
$ws.Cells.Item(29,3).Value = 35
$ws.Cells.Item(29,4).Value = 35
$ws.Cells.Item(29,5).Value = 0
$ws.Cells.Item(29,6).Value = 1
$ws.Cells.Item(29,7).Value = $false
$ws.Cells.Item(29,8).Value = 213

$ws.Cells.Item(30,1).Value = "price"
$ws.Cells.Item(30,2).Value = 0.1379310344827586
$ws.Cells.Item(30,3).Value = 48
$ws.Cells.Item(30,4).Value = 48
$ws.Cells.Item(30,5).Value = 0
$ws.Cells.Item(30,6).Value = 1
$ws.Cells.Item(30,7).Value = $false
$ws.Cells.Item(30,8).Value = 300

$ws.Cells.Item(31,1).Value = "use"
$ws.Cells.Item(31,2).Value = 0.1068493150684932
$ws.Cells.Item(31,3).Value = 39
$ws.Cells.Item(31,4).Value = 39
$ws.Cells.Item(31,5).Value = 0
$ws.Cells.Item(31,6).Value = 1
$ws.Cells.Item(31,7).Value = $false
$ws.Cells.Item(31,8).Value = 326

$ws.Cells.Item(32,1).Value = "little"
$ws.Cells.Item(32,2).Value = 0.0738255033557047
$ws.Cells.Item(32,3).Value = 33
$ws.Cells.Item(32,4).Value = 35
$ws.Cells.Item(32,5).Value = 0.06
$ws.Cells.Item(32,6).Value = 0.94
$ws.Cells.Item(32,7).Value = $true
$ws.Cells.Item(32,8).Value = 414

$ws.Cells.Item(33,1).Value = "much"
$ws.Cells.Item(33,2).Value = 0.06960556844547564
$ws.Cells.Item(33,3).Value = 30
$ws.Cells.Item(33,4).Value = 33
$ws.Cells.Item(33,5).Value = 0.09
$ws.Cells.Item(33,6).Value = 0.91
$ws.Cells.Item(33,7).Value = $true
$ws.Cells.Item(33,8).Value = 401

$ws.Cells.Item(34,1).Value = "like"
$ws.Cells.Item(34,2).Value = 0.05766062602965404
$ws.Cells.Item(34,3).Value = 35
$ws.Cells.Item(34,4).Value = 36
$ws.Cells.Item(34,5).Value = 0.03
$ws.Cells.Item(34,6).Value = 0.97
$ws.Cells.Item(34,7).Value = $true
$ws.Cells.Item(34,8).Value = 572

$ws.Cells.Item(35,1).Value = "toy"
$ws.Cells.Item(35,2).Value = 0.04447852760736196
$ws.Cells.Item(35,3).Value = 29
$ws.Cells.Item(35,4).Value = 32
$ws.Cells.Item(35,5).Value = 0.09
$ws.Cells.Item(35,6).Value = 0.91
$ws.Cells.Item(35,7).Value = $true
$ws.Cells.Item(35,8).Value = 623

$ws.Cells.Item(3,10).Value = "awesome"
$ws.Cells.Item(3,11).Value = 0.8153846153846154
$ws.Cells.Item(3,12).Value = 53
$ws.Cells.Item(3,13).Value = 53
$ws.Cells.Item(3,14).Value = 1
$ws.Cells.Item(3,15).Value = 0
$ws.Cells.Item(3,16).Value = $false
$ws.Cells.Item(3,17).Value = 12

$ws.Cells.Item(4,10).Value = "wonderful"
$ws.Cells.Item(4,11).Value = 0.8035714285714286
$ws.Cells.Item(4,12).Value = 45
$ws.Cells.Item(4,13).Value = 45
$ws.Cells.Item(4,14).Value = 1
$ws.Cells.Item(4,15).Value = 0
$ws.Cells.Item(4,16).Value = $false
$ws.Cells.Item(4,17).Value = 11

$ws.Cells.Item(5,10).Value = "favorite"
$ws.Cells.Item(5,11).Value = 0.6881720430107527
$ws.Cells.Item(5,12).Value = 64
$ws.Cells.Item(5,13).Value = 64
$ws.Cells.Item(5,14).Value = 1
$ws.Cells.Item(5,15).Value = 0
$ws.Cells.Item(5,16).Value = $false
$ws.Cells.Item(5,17).Value = 29

$ws.Cells.Item(6,10).Value = "classic"
$ws.Cells.Item(6,11).Value = 0.5849056603773585
$ws.Cells.Item(6,12).Value = 31
$ws.Cells.Item(6,13).Value = 31
$ws.Cells.Item(6,14).Value = 1
$ws.Cells.Item(6,15).Value = 0
$ws.Cells.Item(6,16).Value = $false
$ws.Cells.Item(6,17).Value = 22

$ws.Cells.Item(7,10).Value = "excellent"
$ws.Cells.Item(7,11).Value = 0.578125
$ws.Cells.Item(7,12).Value = 37
$ws.Cells.Item(7,13).Value = 37
$ws.Cells.Item(7,14).Value = 1
$ws.Cells.Item(7,15).Value = 0
$ws.Cells.Item(7,16).Value = $false
$ws.Cells.Item(7,17).Value = 27

$ws.Cells.Item(8,10).Value = "thank"
$ws.Cells.Item(8,11).Value = 0.4202898550724637
$ws.Cells.Item(8,12).Value = 29
$ws.Cells.Item(8,13).Value = 29
$ws.Cells.Item(8,14).Value = 1
$ws.Cells.Item(8,15).Value = 0
$ws.Cells.Item(8,16).Value = $false
$ws.Cells.Item(8,17).Value = 40

$ws.Cells.Item(9,10).Value = "great"
$ws.Cells.Item(9,11).Value = 0.3860655737704918
$ws.Cells.Item(9,12).Value = 471
$ws.Cells.Item(9,13).Value = 471
$ws.Cells.Item(9,14).Value = 1
$ws.Cells.Item(9,15).Value = 0
$ws.Cells.Item(9,16).Value = $false
$ws.Cells.Item(9,17).Value = 749

$ws.Cells.Item(10,10).Value = "love"
$ws.Cells.Item(10,11).Value = 0.3314203730272597
$ws.Cells.Item(10,12).Value = 231
$ws.Cells.Item(10,13).Value = 231
$ws.Cells.Item(10,14).Value = 1
$ws.Cells.Item(10,15).Value = 0
$ws.Cells.Item(10,16).Value = $false
$ws.Cells.Item(10,17).Value = 466

$ws.Cells.Item(11,10).Value = "loves"
$ws.Cells.Item(11,11).Value = 0.3132780082987552
$ws.Cells.Item(11,12).Value = 151
$ws.Cells.Item(11,13).Value = 151
$ws.Cells.Item(11,14).Value = 1
$ws.Cells.Item(11,15).Value = 0
$ws.Cells.Item(11,16).Value = $false
$ws.Cells.Item(11,17).Value = 331

$ws.Cells.Item(12,10).Value = "perfect"
$ws.Cells.Item(12,11).Value = 0.2349397590361446
$ws.Cells.Item(12,12).Value = 39
$ws.Cells.Item(12,13).Value = 39
$ws.Cells.Item(12,14).Value = 1
$ws.Cells.Item(12,15).Value = 0
$ws.Cells.Item(12,16).Value = $false
$ws.Cells.Item(12,17).Value = 127

$ws.Cells.Item(13,10).Value = "loved"
$ws.Cells.Item(13,11).Value = 0.1926605504587156
$ws.Cells.Item(13,12).Value = 63
$ws.Cells.Item(13,13).Value = 63
$ws.Cells.Item(13,14).Value = 1
$ws.Cells.Item(13,15).Value = 0
$ws.Cells.Item(13,16).Value = $false
$ws.Cells.Item(13,17).Value = 264

$ws.Cells.Item(14,10).Value = "friends"
$ws.Cells.Item(14,11).Value = 0.1587301587301587
$ws.Cells.Item(14,12).Value = 30
$ws.Cells.Item(14,13).Value = 30
$ws.Cells.Item(14,14).Value = 1
$ws.Cells.Item(14,15).Value = 0
$ws.Cells.Item(14,16).Value = $false
$ws.Cells.Item(14,17).Value = 159

$ws.Cells.Item(15,10).Value = "fun"
$ws.Cells.Item(15,11).Value = 0.09465381244522349
$ws.Cells.Item(15,12).Value = 108
$ws.Cells.Item(15,13).Value = 108
$ws.Cells.Item(15,14).Value = 1
$ws.Cells.Item(15,15).Value = 0
$ws.Cells.Item(15,16).Value = $false
$ws.Cells.Item(15,17).Value = 1033

$ws.Cells.Item(16,10).Value = "game"
$ws.Cells.Item(16,11).Value = 0.04421326397919376
$ws.Cells.Item(16,12).Value = 68
$ws.Cells.Item(16,13).Value = 71
$ws.Cells.Item(16,14).Value = 0.96
$ws.Cells.Item(16,15).Value = 0.04000000000000004
$ws.Cells.Item(16,16).Value = $true
$ws.Cells.Item(16,17).Value = 1470

$ws.Range("A34").Copy()
$ws.Range("A35").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("A1:Q35").EntireRow.AutoFit() | Out-Null
